$wb = $excel.ActiveWorkbook

# Overview sheet: file de8179a3-...md is now "Ready for handoff" in both languages
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: status -> "Ready for handoff", new handoff datetime recorded
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-09 04:58:49"

# de-de sheet: status -> "Ready for handoff", new handoff datetime recorded
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-09 04:58:51"
